$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Table19"

$ws.Range("A1").Value = "employeenumber"
$ws.Range("B1").Value = "firstname"
$ws.Range("C1").Value = "lastname"

$ws.Columns.Item(1).ColumnWidth = 17.85546875
$ws.Columns.Item(2).ColumnWidth = 13.7109375
$ws.Columns.Item(3).ColumnWidth = 11.5703125

$ws.Range("B11").Select() | Out-Null
